# Add a new "2022-Q3" sheet right after "总计", pushing the existing
# quarterly sheets one tab to the right, and populate it with the new
# quarter's fund-holding data. Also insert the corresponding summary row
# on the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row for 2022-Q3
#    right after the header row, pushing all existing data rows down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 3
$summary.Cells.Item(2, 4).Value = 0.34

$summary.Cells.Item(2, 1).Font.Bold = $true
$summary.Cells.Item(2, 1).HorizontalAlignment = -4108
$summary.Cells.Item(2, 1).VerticalAlignment = -4160

# Renumber the index column (A) for the rows that followed, since they
# each shifted down by one position (values are simply 0..7 top to bottom).
for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计" and fill it
#    with the fund-holding breakdown for that quarter.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("B1:H1").HorizontalAlignment = -4108
$q3.Range("B1:H1").VerticalAlignment = -4160

$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "009562"
$q3.Cells.Item(2, 3).Value = "工银全球股票（QDII）美元"
$q3.Cells.Item(2, 4).Value = "5.89"
$q3.Cells.Item(2, 5).Value = "93.72"
$q3.Cells.Item(2, 6).Value = "1.90"
$q3.Cells.Item(2, 7).Value = "0.1119"
$q3.Cells.Item(2, 8).Value = 10

$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "009563"
$q3.Cells.Item(3, 3).Value = "工银全球股票（QDII）港币"
$q3.Cells.Item(3, 4).Value = "5.89"
$q3.Cells.Item(3, 5).Value = "93.72"
$q3.Cells.Item(3, 6).Value = "1.90"
$q3.Cells.Item(3, 7).Value = "0.1119"
$q3.Cells.Item(3, 8).Value = 10

$q3.Cells.Item(4, 1).Value = 2
$q3.Cells.Item(4, 2).Value = "486001"
$q3.Cells.Item(4, 3).Value = "工银瑞信中国机会全球配置股票（QDII）人民币"
$q3.Cells.Item(4, 4).Value = "5.89"
$q3.Cells.Item(4, 5).Value = "93.72"
$q3.Cells.Item(4, 6).Value = "1.90"
$q3.Cells.Item(4, 7).Value = "0.1119"
$q3.Cells.Item(4, 8).Value = 10

$q3.Range("A2:A4").Font.Bold = $true
$q3.Range("A2:A4").HorizontalAlignment = -4108
$q3.Range("A2:A4").VerticalAlignment = -4160

$summary.Select()
